$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98: The Dotted Line
$ws.Range("H98").Value = 3809
$ws.Range("I98").Value = 2365.5312
$ws.Range("K98").Value = 2365.5312
$ws.Range("M98").Value = -867.5311999999999

# Row 113: Amaro Kart
$ws.Range("H113").Value = 3045.625
$ws.Range("I113").Value = 2465.7144
$ws.Range("J113").Value = 3496.6667
$ws.Range("K113").Value = 2465.7144
$ws.Range("L113").Value = 3496.6667
$ws.Range("M113").Value = 788.2856000000002
$ws.Range("N113").Value = -10004.6667

# Row 122: Wishful Inking
$ws.Range("H122").Value = 3809
$ws.Range("I122").Value = 2365.5312
$ws.Range("K122").Value = 7096.5936
$ws.Range("M122").Value = -4646.5936

$ws = $wb.Worksheets.Item("ARM")
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1474.875
$ws.Range("I110").Value = 1500
$ws.Range("J110").Value = 1433
$ws.Range("K110").Value = 1500
$ws.Range("L110").Value = 1433
$ws.Range("M110").Value = 545
$ws.Range("N110").Value = -5523

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 3508.9033
$ws.Range("I122").Value = 3359.0386
$ws.Range("J122").Value = 4288.2
$ws.Range("K122").Value = 10077.1158
$ws.Range("L122").Value = 12864.6
$ws.Range("M122").Value = -7627.1158
$ws.Range("N122").Value = -17764.6

# Row 123: The Armoire Is Open
$ws.Range("H123").Value = 24299.166
$ws.Range("J123").Value = 24299.166
$ws.Range("L123").Value = 24299.166
$ws.Range("N123").Value = -34099.166

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 5503.121
$ws.Range("I132").Value = 7565.1055
$ws.Range("J132").Value = 2704.7144
$ws.Range("K132").Value = 22695.3165
$ws.Range("L132").Value = 8114.1432
$ws.Range("M132").Value = -20165.3165
$ws.Range("N132").Value = -13174.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2676.0967
$ws.Range("I134").Value = 2255.3809
$ws.Range("J134").Value = 3559.6
$ws.Range("K134").Value = 6766.1427
$ws.Range("L134").Value = 10678.8
$ws.Range("M134").Value = -4231.1427
$ws.Range("N134").Value = -15748.8

$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 9998.333000000001
$ws.Range("J4").Value = 9998.333000000001
$ws.Range("L4").Value = 9998.333000000001
$ws.Range("N4").Value = -10222.333

# Row 31: Wall Not Found
$ws.Range("H31").Value = 2787.0588
$ws.Range("I31").Value = 1455.7142
$ws.Range("K31").Value = 1455.7142
$ws.Range("M31").Value = -1160.7142

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2787.0588
$ws.Range("I34").Value = 1455.7142
$ws.Range("K34").Value = 1455.7142
$ws.Range("M34").Value = -1253.7142

# Row 97: Wood That You Could
$ws.Range("H97").Value = 31998.5
$ws.Range("J97").Value = 31998.5
$ws.Range("L97").Value = 31998.5
$ws.Range("N97").Value = -33980.5

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 301815.16
$ws.Range("I132").Value = 437017.7
$ws.Range("K132").Value = 1311053.1
$ws.Range("M132").Value = -1308523.1

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 623377.5600000001
$ws.Range("I12").Value = 88.2
$ws.Range("J12").Value = 743240.9399999999
$ws.Range("K12").Value = 264.6
$ws.Range("L12").Value = 2229722.82
$ws.Range("M12").Value = -91.60000000000002
$ws.Range("N12").Value = -2230068.82

# Row 22: A Total Nut Job
$ws.Range("H22").Value = 35716140
$ws.Range("I22").Value = 250000000
$ws.Range("J22").Value = 2163.3333
$ws.Range("K22").Value = 750000000
$ws.Range("L22").Value = 6489.999899999999
$ws.Range("M22").Value = -749999831
$ws.Range("N22").Value = -6827.999899999999

# Row 27: Brain Food
$ws.Range("H27").Value = 35716140
$ws.Range("I27").Value = 250000000
$ws.Range("J27").Value = 2163.3333
$ws.Range("K27").Value = 750000000
$ws.Range("L27").Value = 6489.999899999999
$ws.Range("M27").Value = -749999898
$ws.Range("N27").Value = -6693.999899999999

# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 2264.75
$ws.Range("J39").Value = 2753
$ws.Range("L39").Value = 8259
$ws.Range("N39").Value = -8847

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 3551.3076
$ws.Range("J55").Value = 3551.3076
$ws.Range("L55").Value = 10653.9228
$ws.Range("N55").Value = -11007.9228

# Row 76: Old Victories, New Tastes
$ws.Range("H76").Value = 2353.5
$ws.Range("I76").Value = 1138
$ws.Range("J76").Value = 6000
$ws.Range("K76").Value = 3414
$ws.Range("L76").Value = 18000
$ws.Range("M76").Value = -3031
$ws.Range("N76").Value = -18766

# Row 79: The Eats of Authenticity (L)
$ws.Range("H79").Value = 2353.5
$ws.Range("I79").Value = 1138
$ws.Range("J79").Value = 6000
$ws.Range("K79").Value = 3414
$ws.Range("L79").Value = 18000
$ws.Range("M79").Value = -2088
$ws.Range("N79").Value = -20652

# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 15250
$ws.Range("I80").Value = 15000
$ws.Range("J80").Value = 15400
$ws.Range("K80").Value = 45000
$ws.Range("L80").Value = 46200
$ws.Range("M80").Value = -44064
$ws.Range("N80").Value = -48072

# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 15250
$ws.Range("I83").Value = 15000
$ws.Range("J83").Value = 15400
$ws.Range("K83").Value = 135000
$ws.Range("L83").Value = 138600
$ws.Range("M83").Value = -130320
$ws.Range("N83").Value = -147960

# Row 100: Souper
$ws.Range("H100").Value = 11658
$ws.Range("J100").Value = 2263.3333
$ws.Range("L100").Value = 6789.999899999999
$ws.Range("N100").Value = -8411.999899999999

# Row 103: West Meats East
$ws.Range("H103").Value = 3974.5789
$ws.Range("I103").Value = 2144.2
$ws.Range("J103").Value = 4628.2856
$ws.Range("K103").Value = 6432.599999999999
$ws.Range("L103").Value = 13884.8568
$ws.Range("M103").Value = -5553.599999999999
$ws.Range("N103").Value = -15642.8568

# Row 106: Herky Jerky
$ws.Range("H106").Value = 5562.857
$ws.Range("J106").Value = 5562.857
$ws.Range("L106").Value = 16688.571
$ws.Range("N106").Value = -18580.571

# Row 109: Cure for What Ails
$ws.Range("H109").Value = 3109.3333
$ws.Range("J109").Value = 4395
$ws.Range("L109").Value = 13185
$ws.Range("N109").Value = -15265

# Row 112: Sweet Tooth
$ws.Range("H112").Value = 4384.722
$ws.Range("I112").Value = 775
$ws.Range("J112").Value = 4835.9375
$ws.Range("K112").Value = 2325
$ws.Range("L112").Value = 14507.8125
$ws.Range("M112").Value = -1217
$ws.Range("N112").Value = -16723.8125

# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 4808.353
$ws.Range("I136").Value = 952.8570999999999
$ws.Range("J136").Value = 7507.2
$ws.Range("K136").Value = 2858.5713
$ws.Range("L136").Value = 22521.6
$ws.Range("M136").Value = 2241.4287
$ws.Range("N136").Value = -32721.6

$ws = $wb.Worksheets.Item("GSM")
# Row 93: One Ring Circus
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 109: You're My Wonderhall
$ws.Range("H109").Value = 9159.058999999999
$ws.Range("J109").Value = 9159.058999999999
$ws.Range("L109").Value = 9159.058999999999
$ws.Range("N109").Value = -11239.059

# Row 123: Workplace Workout
$ws.Range("H123").Value = 8619
$ws.Range("J123").Value = 8619
$ws.Range("L123").Value = 8619
$ws.Range("N123").Value = -13519

# Row 132: On Board for Lar
$ws.Range("H132").Value = 2821.0264
$ws.Range("I132").Value = 2369.889
$ws.Range("K132").Value = 7109.667
$ws.Range("M132").Value = -4579.667

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1378.4286
$ws.Range("I22").Value = 1530
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 1530
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -1235
$ws.Range("N22").Value = -1589.5

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1378.4286
$ws.Range("I27").Value = 1530
$ws.Range("J27").Value = 999.5
$ws.Range("K27").Value = 1530
$ws.Range("L27").Value = 999.5
$ws.Range("M27").Value = -1423
$ws.Range("N27").Value = -1213.5

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 2797.5535
$ws.Range("I132").Value = 1900.4517
$ws.Range("J132").Value = 3909.96
$ws.Range("K132").Value = 5701.355100000001
$ws.Range("L132").Value = 11729.88
$ws.Range("M132").Value = -3171.355100000001
$ws.Range("N132").Value = -16789.88

$ws = $wb.Worksheets.Item("WVR")
# Row 11: Wiggle Room
$ws.Range("H11").Value = 50000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Row 107: Flax Wax
$ws.Range("H107").Value = 417.125
$ws.Range("I107").Value = 367.33334
$ws.Range("J107").Value = 447
$ws.Range("K107").Value = 1102.00002
$ws.Range("L107").Value = 1341
$ws.Range("M107").Value = 817.9999800000001
$ws.Range("N107").Value = -5181

# Row 123: Helping Handwear
$ws.Range("H123").Value = 23787.773
$ws.Range("J123").Value = 23787.773
$ws.Range("L123").Value = 23787.773
$ws.Range("N123").Value = -33587.773

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1965.0222
$ws.Range("I132").Value = 1169.7188
$ws.Range("J132").Value = 3922.6924
$ws.Range("K132").Value = 3509.1564
$ws.Range("L132").Value = 11768.0772
$ws.Range("M132").Value = -979.1564000000003
$ws.Range("N132").Value = -16828.0772
